$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A80").Value = "Biziura lobata"
$ws.Range("B80").Value = "https://www.xeno-canto.org/201462/download"
$ws.Range("C80").Value = "Dunn's Swamp, New South Wales, Australia"
$ws.Range("D80").Value = "Marc Anderson"
$ws.Range("E80").Value = "https://www.xeno-canto.org/contributor/EHGWCIGILC"
$ws.Range("F80").Value = "assets/misc/cc.png"
$ws.Range("G80").Value = "https://creativecommons.org/licenses/by-nc-nd/4.0/"
$ws.Range("A81").Value = "Bubulcus ibis"
$ws.Range("B81").Value = "https://www.xeno-canto.org/393573/download"
$ws.Range("C81").Value = "Maleny, Queensland, Australia"
$ws.Range("D81").Value = "Greg McLachlan"
$ws.Range("E81").Value = "https://www.xeno-canto.org/contributor/MXMFTGKZDR"
$ws.Range("F81").Value = "assets/misc/cc.png"
$ws.Range("G81").Value = "https://creativecommons.org/licenses/by-nc-sa/4.0/"
$ws.Range("A82").Value = "Burhinus grallarius"
$ws.Range("B82").Value = "https://www.xeno-canto.org/439120/download"
$ws.Range("C82").Value = "Coen River, Queensland, Australia"
$ws.Range("D82").Value = "Marc Anderson"
$ws.Range("E82").Value = "https://www.xeno-canto.org/contributor/EHGWCIGILC"
$ws.Range("F82").Value = "assets/misc/cc.png"
$ws.Range("G82").Value = "https://creativecommons.org/licenses/by-nc-nd/4.0/"
$ws.Range("A83").Value = "Burhinus grallarius"
$ws.Range("B83").Value = "https://www.xeno-canto.org/382911/download"
$ws.Range("C83").Value = "Dryandra Woodland, Western Australia, Australia"
$ws.Range("D83").Value = "Marc Anderson"
$ws.Range("E83").Value = "https://www.xeno-canto.org/contributor/EHGWCIGILC"
$ws.Range("F83").Value = "assets/misc/cc.png"
$ws.Range("G83").Value = "https://creativecommons.org/licenses/by-nc-nd/4.0/"
$ws.Range("A84").Value = "Burhinus grallarius"
$ws.Range("B84").Value = "https://www.xeno-canto.org/174119/download"
$ws.Range("C84").Value = "Stewart Creek Valley, Queensland, Australia"
$ws.Range("D84").Value = "Marc Anderson"
$ws.Range("E84").Value = "https://www.xeno-canto.org/contributor/EHGWCIGILC"
$ws.Range("F84").Value = "assets/misc/cc.png"
$ws.Range("G84").Value = "https://creativecommons.org/licenses/by-nc-nd/4.0/"
$ws.Range("A85").Value = "Butorides striata"
$ws.Range("B85").Value = "https://www.xeno-canto.org/359600/download"
$ws.Range("C85").Value = "Sabah, Borneo, Malaysia"
$ws.Range("D85").Value = "Peter Boesman"
$ws.Range("E85").Value = "https://www.xeno-canto.org/contributor/OOECIWCSWV"
$ws.Range("F85").Value = "assets/misc/cc.png"
$ws.Range("G85").Value = "https://creativecommons.org/licenses/by-nc-nd/4.0/"
$ws.Range("A86").Value = "Butorides striata"
$ws.Range("B86").Value = "https://www.xeno-canto.org/210507/download"
$ws.Range("C86").Value = "Tweed Heads, New South Wales, Australia"
$ws.Range("D86").Value = "Nick Talbot"
$ws.Range("E86").Value = "https://www.xeno-canto.org/contributor/CCUCXWCPSW"
$ws.Range("F86").Value = "assets/misc/cc.png"
$ws.Range("G86").Value = "https://creativecommons.org/licenses/by-nc-sa/4.0/"
$ws.Range("A87").Value = "Cacatua galerita"
$ws.Range("B87").Value = "https://www.xeno-canto.org/438828/download"
$ws.Range("C87").Value = "Oyala-Thumotang National Park, Queensland, Australia"
$ws.Range("D87").Value = "Marc Anderson"
$ws.Range("E87").Value = "https://www.xeno-canto.org/contributor/EHGWCIGILC"
$ws.Range("F87").Value = "assets/misc/cc.png"
$ws.Range("G87").Value = "https://creativecommons.org/licenses/by-nc-nd/4.0/"
$ws.Range("A88").Value = "Cacatua galerita"
$ws.Range("B88").Value = "https://www.xeno-canto.org/380456/download"
$ws.Range("C88").Value = "Long Point, New South Wales, Australia"
$ws.Range("D88").Value = "Greg McLachlan"
$ws.Range("E88").Value = "https://www.xeno-canto.org/contributor/MXMFTGKZDR"
$ws.Range("F88").Value = "assets/misc/cc.png"
$ws.Range("G88").Value = "https://creativecommons.org/licenses/by-nc-sa/4.0/"
$ws.Range("A89").Value = "Cacatua sanguinea"
$ws.Range("B89").Value = "https://www.xeno-canto.org/407697/download"
$ws.Range("C89").Value = "Lake Hattah, Victoria, Australia"
$ws.Range("D89").Value = "Frank Lambert"
$ws.Range("E89").Value = "https://www.xeno-canto.org/contributor/YTUXOCTUEM"
$ws.Range("F89").Value = "assets/misc/cc.png"
$ws.Range("G89").Value = "https://creativecommons.org/licenses/by-nc-nd/4.0/"
$ws.Range("A90").Value = "Cacatua sanguinea"
$ws.Range("B90").Value = "https://www.xeno-canto.org/320990/download"
$ws.Range("C90").Value = "Maffra, Victoria, Australia"
$ws.Range("D90").Value = "Nick Talbot"
$ws.Range("E90").Value = "https://www.xeno-canto.org/contributor/CCUCXWCPSW"
$ws.Range("F90").Value = "assets/misc/cc.png"
$ws.Range("G90").Value = "https://creativecommons.org/licenses/by-nc-sa/4.0/"
$ws.Range("A91").Value = "Cacatua tenuirostris"
$ws.Range("B91").Value = "https://www.xeno-canto.org/200313/download"
$ws.Range("C91").Value = "Wyperfeld National Park, Victoria, Australia"
$ws.Range("D91").Value = "Marc Anderson"
$ws.Range("E91").Value = "https://www.xeno-canto.org/contributor/EHGWCIGILC"
$ws.Range("F91").Value = "assets/misc/cc.png"
$ws.Range("G91").Value = "https://creativecommons.org/licenses/by-nc-nd/4.0/"
$ws.Range("A92").Value = "Cacatua tenuirostris"
$ws.Range("B92").Value = "https://www.xeno-canto.org/187680/download"
$ws.Range("C92").Value = "Adelaide, South Australia, Australia"
$ws.Range("D92").Value = "Nick Talbot"
$ws.Range("E92").Value = "https://www.xeno-canto.org/contributor/CCUCXWCPSW"
$ws.Range("F92").Value = "assets/misc/cc.png"
$ws.Range("G92").Value = "https://creativecommons.org/licenses/by-nc-sa/4.0/"
$ws.Range("A93").Value = "Cacomantis flabelliformis"
$ws.Range("B93").Value = "https://www.xeno-canto.org/407773/download"
$ws.Range("C93").Value = "South Bruny Island, Tasmania, Australia"
$ws.Range("D93").Value = "Frank Lambert"
$ws.Range("E93").Value = "https://www.xeno-canto.org/contributor/YTUXOCTUEM"
$ws.Range("F93").Value = "assets/misc/cc.png"
$ws.Range("G93").Value = "https://creativecommons.org/licenses/by-nc-nd/4.0/"
$ws.Range("A94").Value = "Cacomantis flabelliformis"
$ws.Range("B94").Value = "https://www.xeno-canto.org/382702/download"
$ws.Range("C94").Value = "Cheynes Beach, Western Australia, Australia"
$ws.Range("D94").Value = "Marc Anderson"
$ws.Range("E94").Value = "https://www.xeno-canto.org/contributor/EHGWCIGILC"
$ws.Range("F94").Value = "assets/misc/cc.png"
$ws.Range("G94").Value = "https://creativecommons.org/licenses/by-nc-nd/4.0/"
$ws.Range("A95").Value = "Cacomantis flabelliformis"
$ws.Range("B95").Value = "https://www.xeno-canto.org/140206/download"
$ws.Range("C95").Value = "Barren Grounds, New South Wales, Australia"
$ws.Range("D95").Value = "Marc Anderson"
$ws.Range("E95").Value = "https://www.xeno-canto.org/contributor/EHGWCIGILC"
$ws.Range("F95").Value = "assets/misc/cc.png"
$ws.Range("G95").Value = "https://creativecommons.org/licenses/by-nc-nd/4.0/"
$ws.Range("A96").Value = "Cacomantis pallidus"
$ws.Range("B96").Value = "https://www.xeno-canto.org/389397/download"
$ws.Range("C96").Value = "Pitt Town Lagoon, New South Wales, Australia"
$ws.Range("D96").Value = "Marc Anderson"
$ws.Range("E96").Value = "https://www.xeno-canto.org/contributor/EHGWCIGILC"
$ws.Range("F96").Value = "assets/misc/cc.png"
$ws.Range("G96").Value = "https://creativecommons.org/licenses/by-nc-nd/4.0/"
$ws.Range("A97").Value = "Cacomantis pallidus"
$ws.Range("B97").Value = "https://www.xeno-canto.org/389393/download"
$ws.Range("C97").Value = "Pitt Town Lagoon, New South Wales, Australia"
$ws.Range("D97").Value = "Marc Anderson"
$ws.Range("E97").Value = "https://www.xeno-canto.org/contributor/EHGWCIGILC"
$ws.Range("F97").Value = "assets/misc/cc.png"
$ws.Range("G97").Value = "https://creativecommons.org/licenses/by-nc-nd/4.0/"
$ws.Range("A98").Value = "Cacomantis variolosus"
$ws.Range("B98").Value = "https://www.xeno-canto.org/287384/download"
$ws.Range("C98").Value = "Barcoongere State Forest, New South Wales, Australia"
$ws.Range("D98").Value = "Marc Anderson"
$ws.Range("E98").Value = "https://www.xeno-canto.org/contributor/EHGWCIGILC"
$ws.Range("F98").Value = "assets/misc/cc.png"
$ws.Range("G98").Value = "https://creativecommons.org/licenses/by-nc-nd/4.0/"
$ws.Range("A99").Value = "Cacomantis variolosus"
$ws.Range("B99").Value = "https://www.xeno-canto.org/201458/download"
$ws.Range("C99").Value = "Dunn's Swamp, New South Wales, Australia"
$ws.Range("D99").Value = "Marc Anderson"
$ws.Range("E99").Value = "https://www.xeno-canto.org/contributor/EHGWCIGILC"
$ws.Range("F99").Value = "assets/misc/cc.png"
$ws.Range("G99").Value = "https://creativecommons.org/licenses/by-nc-nd/4.0/"
$ws.Range("A100").Value = "Calidris acuminata"
$ws.Range("B100").Value = "https://www.xeno-canto.org/107605/download"
$ws.Range("C100").Value = "St Paul Island, Alaska, United States"
$ws.Range("D100").Value = "Ryan O'Donnell"
$ws.Range("E100").Value = "https://www.xeno-canto.org/contributor/SDXVTLDNGJ"
$ws.Range("F100").Value = "assets/misc/cc.png"
$ws.Range("G100").Value = "https://creativecommons.org/licenses/by-nc-nd/2.5/"
$ws.Range("A101").Value = "Calidris acuminata"
$ws.Range("B101").Value = "https://www.xeno-canto.org/283313/download"
$ws.Range("C101").Value = "Tuggerah, New South Wales, Australia"
$ws.Range("D101").Value = "Marc Anderson"
$ws.Range("E101").Value = "https://www.xeno-canto.org/contributor/EHGWCIGILC"
$ws.Range("F101").Value = "assets/misc/cc.png"
$ws.Range("G101").Value = "https://creativecommons.org/licenses/by-nc-nd/4.0/"
$ws.Range("A102").Value = "Calidris alba"
$ws.Range("B102").Value = "https://www.xeno-canto.org/412103/download"
$ws.Range("C102").Value = "Orog Lake, Mongolia"
$ws.Range("D102").Value = "Frank Lambert"
$ws.Range("E102").Value = "https://www.xeno-canto.org/contributor/YTUXOCTUEM"
$ws.Range("F102").Value = "assets/misc/cc.png"
$ws.Range("G102").Value = "https://creativecommons.org/licenses/by-nc-nd/4.0/"
$ws.Range("A103").Value = "Calidris alba"
$ws.Range("B103").Value = "https://www.xeno-canto.org/223302/download"
$ws.Range("C103").Value = "Tarifa, Andalusia, Spain"
$ws.Range("D103").Value = "Karri Kuitunen"
$ws.Range("E103").Value = "https://www.xeno-canto.org/contributor/XJIOOFMPPX"
$ws.Range("F103").Value = "assets/misc/cc.png"
$ws.Range("G103").Value = "https://creativecommons.org/licenses/by-nc-sa/4.0/"
$ws.Range("A104").Value = "Calidris canutus"
$ws.Range("B104").Value = "https://www.xeno-canto.org/416902/download"
$ws.Range("C104").Value = "Balsnes, Troms, Norway"
$ws.Range("E104").Value = "https://www.xeno-canto.org/contributor/TLPLNAINFU"
$ws.Range("D104").Value = "Stein Ø. Nilsen"
$ws.Range("F104").Value = "assets/misc/cc.png"
$ws.Range("G104").Value = "https://creativecommons.org/licenses/by-nc-sa/4.0/"
$ws.Range("A105").Value = "Calidris falcinellus"
$ws.Range("B105").Value = "https://www.xeno-canto.org/342360/download"
$ws.Range("C105").Value = "Iisakkiaapa, Lapland, Finland"
$ws.Range("D105").Value = "Tero Linjama"
$ws.Range("E105").Value = "https://www.xeno-canto.org/contributor/YSDNMROVID"
$ws.Range("F105").Value = "assets/misc/cc.png"
$ws.Range("G105").Value = "https://creativecommons.org/licenses/by-nc-sa/4.0/"
$ws.Range("A106").Value = "Calidris ferruginea"
$ws.Range("B106").Value = "https://www.xeno-canto.org/346641/download"
$ws.Range("C106").Value = "Walvisbay, Namibia"
$ws.Range("D106").Value = "Peter Boesman"
$ws.Range("E106").Value = "https://www.xeno-canto.org/contributor/OOECIWCSWV"
$ws.Range("F106").Value = "assets/misc/cc.png"
$ws.Range("G106").Value = "https://creativecommons.org/licenses/by-nc-nd/4.0/"
$ws.Range("A107").Value = "Calidris ferruginea"
$ws.Range("B107").Value = "https://www.xeno-canto.org/184321/download"
$ws.Range("C107").Value = "Kazaly District, Kyzylorda Province, Kazakhstan"
$ws.Range("D107").Value = "Albert Lastukhin"
$ws.Range("E107").Value = "https://www.xeno-canto.org/contributor/LELYWQKUZX"
$ws.Range("F107").Value = "assets/misc/cc.png"
$ws.Range("G107").Value = "https://creativecommons.org/licenses/by-nc-sa/4.0/"
$ws.Range("A108").Value = "Calidris melanotos"
$ws.Range("B108").Value = "https://www.xeno-canto.org/406470/download"
$ws.Range("C108").Value = "Barrow, Alaksa, United States"
$ws.Range("D108").Value = "Patrik Åberg"
$ws.Range("E108").Value = "https://www.xeno-canto.org/contributor/BPSDQEOJWG"
$ws.Range("F108").Value = "assets/misc/cc.png"
$ws.Range("G108").Value = "https://creativecommons.org/licenses/by-nc-sa/4.0/"
$ws.Range("A109").Value = "Calidris ruficollis"
$ws.Range("B109").Value = "https://www.xeno-canto.org/107881/download"
$ws.Range("C109").Value = "Monkey Mia, Western Australia, Australia"
$ws.Range("D109").Value = "Matthias Feuersenger"
$ws.Range("E109").Value = "https://www.xeno-canto.org/contributor/HBPYQXTJEV"
$ws.Range("F109").Value = "assets/misc/cc.png"
$ws.Range("G109").Value = "https://creativecommons.org/licenses/by-nc-nd/2.5/"
$ws.Range("A110").Value = "Calidris subminuta"
$ws.Range("B110").Value = "https://www.xeno-canto.org/295966/download"
$ws.Range("C110").Value = "Dalian, China"
$ws.Range("D110").Value = "Tom Beeke"
$ws.Range("E110").Value = "https://www.xeno-canto.org/contributor/CTULRNLZWS"
$ws.Range("F110").Value = "assets/misc/cc.png"
$ws.Range("G110").Value = "https://creativecommons.org/licenses/by-nc-sa/4.0/"
$ws.Range("A111").Value = "Calidris subruficollis"
$ws.Range("B111").Value = "https://www.xeno-canto.org/435877/download"
$ws.Range("C111").Value = "Reserva Natural Palmarí, Rio Javarí, Brazil"
$ws.Range("D111").Value = "Jerome Fischer"
$ws.Range("E111").Value = "https://www.xeno-canto.org/contributor/JPBSNBUUEF"
$ws.Range("F111").Value = "assets/misc/cc.png"
$ws.Range("G111").Value = "https://creativecommons.org/licenses/by-nc-sa/4.0/"
$ws.Range("A112").Value = "Calidris tenuirostris"
$ws.Range("B112").Value = "https://www.xeno-canto.org/396375/download"
$ws.Range("C112").Value = "Trevozhnaya, Chukotka Autonomous Okrug, Russian Federation"
$ws.Range("D112").Value = "Christian A. Jensen"
$ws.Range("E112").Value = "https://www.xeno-canto.org/contributor/IBLQAJNUOV"
$ws.Range("F112").Value = "assets/misc/cc.png"
$ws.Range("G112").Value = "https://creativecommons.org/licenses/by-nc-sa/4.0/"
$ws.Range("A113").Value = "Caligavis chrysops"
$ws.Range("B113").Value = "https://www.xeno-canto.org/390581/download"
$ws.Range("C113").Value = "Burralow Creek, New South Wales, Australia"
$ws.Range("D113").Value = "Marc Anderson"
$ws.Range("E113").Value = "https://www.xeno-canto.org/contributor/EHGWCIGILC"
$ws.Range("F113").Value = "assets/misc/cc.png"
$ws.Range("G113").Value = "https://creativecommons.org/licenses/by-nc-nd/4.0/"
$ws.Range("A114").Value = "Caligavis chrysops"
$ws.Range("B114").Value = "https://www.xeno-canto.org/104876/download"
$ws.Range("C114").Value = "Mogo Creek, New South Wales, Australia"
$ws.Range("D114").Value = "Eliot Miller"
$ws.Range("E114").Value = "https://www.xeno-canto.org/contributor/YLVIJORHMB"
$ws.Range("F114").Value = "assets/misc/cc.png"
$ws.Range("G114").Value = "https://creativecommons.org/licenses/by-nc-nd/2.5/"
$ws.Range("A115").Value = "Calyptorhynchus banksii"
$ws.Range("B115").Value = "https://www.xeno-canto.org/439665/download"
$ws.Range("C115").Value = "Twin Bridges, Queensland, Australia"
$ws.Range("D115").Value = "Marc Anderson"
$ws.Range("E115").Value = "https://www.xeno-canto.org/contributor/EHGWCIGILC"
$ws.Range("F115").Value = "assets/misc/cc.png"
$ws.Range("G115").Value = "https://creativecommons.org/licenses/by-nc-nd/4.0/"
$ws.Range("A116").Value = "Calyptorhynchus banksii"
$ws.Range("B116").Value = "https://www.xeno-canto.org/104983/download"
$ws.Range("C116").Value = "Lakefield National Park, Queensland, Australia"
$ws.Range("D116").Value = "Eliot Miller"
$ws.Range("E116").Value = "https://www.xeno-canto.org/contributor/YLVIJORHMB"
$ws.Range("F116").Value = "assets/misc/cc.png"
$ws.Range("G116").Value = "https://creativecommons.org/licenses/by-nc-nd/2.5/"

$ws.Range("D117").Select()
